$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.937.09"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "3.565.60"
$ws.Range("E3").Value = "  -1.74%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'238.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("D6").Value = "'654.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("E7").Value = "  +10.02%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'1.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.90%  "

$ws.Range("D11").Value = "3.564.23"
$ws.Range("E11").Value = "  -1.77%  "

$ws.Range("D12").Value = "'43.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("E13").Value = "  +0.72%  "

$ws.Range("D14").Value = "'6.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").Value = "4.227.31"
$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("D16").Value = "95.874.11"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").Value = "'0.0000259"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.58%  "

$ws.Range("D18").Value = "3.550.75"
$ws.Range("E18").Value = "  -1.86%  "

$ws.Range("D19").Value = "'7.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.99%  "

$ws.Range("E20").Value = "  -3.50%  "

$ws.Range("E21").Value = "  -1.94%  "

$ws.Range("D22").Value = "'0.516"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.72%  "

$ws.Range("E23").Value = "  -5.94%  "

$ws.Range("D24").Value = "'501.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("D25").Value = "'6.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.27%  "

$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").Value = "'95.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.45%  "

$ws.Range("D28").Value = "'12.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").Value = "3.756.16"
$ws.Range("E29").Value = "  -1.54%  "

$ws.Range("D30").Value = "'0.151"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.00%  "

$ws.Range("E31").Value = "  -4.75%  "

$ws.Range("D32").Value = "'11.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("D34").Value = "'0.182"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.36%  "

$ws.Range("D35").Value = "'0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("D36").Value = "'31.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("D37").Value = "'8.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.21%  "

$ws.Range("D38").Value = "'613.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.99%  "

$ws.Range("D39").Value = "'0.562"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("E40").Value = "  +9.41%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").Value = "'0.900"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.76%  "

$ws.Range("E44").Value = "  +4.87%  "

$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("E46").Value = "  -0.86%  "

$ws.Range("D47").Value = "'0.0421"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("D48").Value = "'2.25"
$ws.Range("D48").Style = "Normal"

$ws.Range("E49").Value = "  -5.53%  "

$ws.Range("E50").Value = "  -0.89%  "

$ws.Range("D51").Value = "'8.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
